# Weekly update: insert this week's two new "Betarraga" price rows at the
# top of the date-ordered block (rows 306:307), pushing the previously
# existing rows down by two (306->308 ... 341->343).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows, shifting rows 306:341 down to 308:343.
$ws.Range("A306:R307").Insert()

# New row 306 - "Primera" quality, week of 2022-07-27.
$ws.Range("A306").Value = 1
$ws.Range("B306").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C306").Value = "Arica y Parinacota"
$ws.Range("D306").Value = 44769
$ws.Range("E306").Value = 15
$ws.Range("F306").Value = 100114014
$ws.Range("G306").Value = "Betarraga"
$ws.Range("H306").Value = "Sin especificar"
$ws.Range("I306").Value = "Primera"
$ws.Range("J306").Value = 800
$ws.Range("K306").Value = 450
$ws.Range("L306").Value = 500
$ws.Range("M306").Value = 475
$ws.Range("N306").Value = "`$/paquete 4 unidades"
$ws.Range("O306").Value = "Región de Arica y Parinacota"
$ws.Range("P306").Value = 119
$ws.Range("Q306").Value = 4
$ws.Range("R306").Value = "Hortaliza"

# New row 307 - "Segunda" quality, week of 2022-07-27.
$ws.Range("A307").Value = 1
$ws.Range("B307").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C307").Value = "Arica y Parinacota"
$ws.Range("D307").Value = 44769
$ws.Range("E307").Value = 15
$ws.Range("F307").Value = 100114014
$ws.Range("G307").Value = "Betarraga"
$ws.Range("H307").Value = "Sin especificar"
$ws.Range("I307").Value = "Segunda"
$ws.Range("J307").Value = 1000
$ws.Range("K307").Value = 450
$ws.Range("L307").Value = 500
$ws.Range("M307").Value = 475
$ws.Range("N307").Value = "`$/paquete 5 unidades"
$ws.Range("O307").Value = "Región de Arica y Parinacota"
$ws.Range("P307").Value = 95
$ws.Range("Q307").Value = 5
$ws.Range("R307").Value = "Hortaliza"
